# Add a new, country-specific "columnsFertilityF1a" parameter row just above
# the existing "columnsFertilityF1b" row on the ColumnsNumberParameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# "columnsFertilityF1b" currently lives on row 27 - insert a new blank row
# above it (shifting it, and everything below, down by one row) and fill it in.
$ws.Rows("27:27").Insert()

$ws.Range("A27").Value = "columnsFertilityF1a"

# Match the existing "quoted number" text formatting used for the other
# numeric-looking KEY/VALUE entries near the top of the sheet (e.g. B3/B4).
$ws.Range("B27").ClearFormats()
$ws.Range("B27").NumberFormat = "0"
$ws.Range("B27").Value = "'5"

# Give column A a bit more breathing room now that labels are a touch longer.
$ws.Range("A:A").ColumnWidth = 32.8

# Restore the originally-selected cell / view.
$ws.Range("B33").Select() | Out-Null
